# Auto-generated Excel COM-interop script to apply the diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 598.4286
$ws.Range("I2").Value = 457.8
$ws.Range("K2").Value = 457.8
$ws.Range("M2").Value = -344.8
$ws.Range("H19").Value = 804.24
$ws.Range("I19").Value = 506.64706
$ws.Range("K19").Value = 506.64706
$ws.Range("M19").Value = -331.64706
$ws.Range("H80").Value = 8017.115
$ws.Range("I80").Value = 6144.2354
$ws.Range("K80").Value = 18432.7062
$ws.Range("M80").Value = -17434.7062
$ws.Range("H83").Value = 8017.115
$ws.Range("I83").Value = 6144.2354
$ws.Range("K83").Value = 55298.11859999999
$ws.Range("M83").Value = -50306.11859999999
$ws.Range("H97").Value = 698.6
$ws.Range("J97").Value = 698.6
$ws.Range("L97").Value = 2095.8
$ws.Range("N97").Value = -3087.8
$ws.Range("H116").Value = 4299.6665
$ws.Range("I116").Value = 3559.4
$ws.Range("J116").Value = 4828.4287
$ws.Range("K116").Value = 3559.4
$ws.Range("L116").Value = 4828.4287
$ws.Range("M116").Value = -117.4000000000001
$ws.Range("N116").Value = -11712.4287
$ws.Range("H129").Value = 2215.3447
$ws.Range("I129").Value = 1661.875
$ws.Range("J129").Value = 2426.1904
$ws.Range("K129").Value = 4985.625
$ws.Range("L129").Value = 7278.5712
$ws.Range("M129").Value = 14.375
$ws.Range("N129").Value = -17278.5712
$ws.Range("H132").Value = 4457.0586
$ws.Range("I132").Value = 4715.893
$ws.Range("K132").Value = 14147.679
$ws.Range("M132").Value = -11617.679
$ws.Range("H141").Value = 6716
$ws.Range("I141").Value = 6716
$ws.Range("K141").Value = 20148
$ws.Range("M141").Value = -14968

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1787.238
$ws.Range("J2").Value = 2169.4
$ws.Range("L2").Value = 2169.4
$ws.Range("N2").Value = -2395.4
$ws.Range("H116").Value = 1787.238
$ws.Range("J116").Value = 2169.4
$ws.Range("L116").Value = 2169.4
$ws.Range("N116").Value = -6757.4
$ws.Range("H122").Value = 8198
$ws.Range("H132").Value = 5732.3335
$ws.Range("J132").Value = 8433.799999999999
$ws.Range("L132").Value = 25301.4
$ws.Range("N132").Value = -30361.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1787.238
$ws.Range("J3").Value = 2169.4
$ws.Range("L3").Value = 2169.4
$ws.Range("N3").Value = -2397.4
$ws.Range("H52").Value = 9000
$ws.Range("I52").Value = 9000
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 9000
$ws.Range("N52").ClearContents()
$ws.Range("M52").Value = -8737
$ws.Range("H121").Value = 9000
$ws.Range("I121").Value = 9000
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 9000
$ws.Range("N121").ClearContents()
$ws.Range("M121").Value = -7253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2476.3076
$ws.Range("I16").Value = 2410.3333
$ws.Range("J16").Value = 2624.75
$ws.Range("K16").Value = 2410.3333
$ws.Range("L16").Value = 2624.75
$ws.Range("M16").Value = -2123.3333
$ws.Range("N16").Value = -3198.75
$ws.Range("H22").Value = 342.72726
$ws.Range("I22").Value = 369.625
$ws.Range("J22").Value = 271
$ws.Range("K22").Value = 369.625
$ws.Range("L22").Value = 271
$ws.Range("M22").Value = -19.625
$ws.Range("N22").Value = -971
$ws.Range("H31").Value = 7631.3335
$ws.Range("I31").Value = 5369.143
$ws.Range("K31").Value = 5369.143
$ws.Range("M31").Value = -5074.143
$ws.Range("H34").Value = 7631.3335
$ws.Range("I34").Value = 5369.143
$ws.Range("K34").Value = 5369.143
$ws.Range("M34").Value = -5167.143
$ws.Range("H58").Value = 6789.7
$ws.Range("I58").Value = 3833.3333
$ws.Range("J58").Value = 8056.7144
$ws.Range("K58").Value = 3833.3333
$ws.Range("L58").Value = 8056.7144
$ws.Range("M58").Value = -3630.3333
$ws.Range("N58").Value = -8462.714400000001
$ws.Range("H86").Value = 4216.4287
$ws.Range("I86").Value = 3626.5881
$ws.Range("K86").Value = 3626.5881
$ws.Range("M86").Value = -2503.5881
$ws.Range("H89").Value = 4216.4287
$ws.Range("I89").Value = 3626.5881
$ws.Range("K89").Value = 18132.9405
$ws.Range("M89").Value = -12516.9405
$ws.Range("H112").Value = 62685.8
$ws.Range("J112").Value = 62685.8
$ws.Range("L112").Value = 62685.8
$ws.Range("N112").Value = -65639.8
$ws.Range("H113").Value = 2476.3076
$ws.Range("I113").Value = 2410.3333
$ws.Range("J113").Value = 2624.75
$ws.Range("K113").Value = 2410.3333
$ws.Range("L113").Value = 2624.75
$ws.Range("M113").Value = -240.3332999999998
$ws.Range("N113").Value = -6964.75
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 6887
$ws.Range("I134").Value = 1542.25
$ws.Range("K134").Value = 4626.75
$ws.Range("M134").Value = -2091.75
$ws.Range("H136").Value = 6789.7
$ws.Range("I136").Value = 3833.3333
$ws.Range("J136").Value = 8056.7144
$ws.Range("K136").Value = 11499.9999
$ws.Range("L136").Value = 24170.1432
$ws.Range("M136").Value = -8949.999899999999
$ws.Range("N136").Value = -29270.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1421.8438
$ws.Range("I2").Value = 447.0909
$ws.Range("J2").Value = 3566.3
$ws.Range("K2").Value = 2682.5454
$ws.Range("L2").Value = 21397.8
$ws.Range("M2").Value = -2569.5454
$ws.Range("N2").Value = -21623.8
$ws.Range("H63").Value = 3725
$ws.Range("I63").Value = 2450
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 7350
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = -6601
$ws.Range("N63").Value = -16498
$ws.Range("H64").Value = 3552.375
$ws.Range("I64").Value = 2104.75
$ws.Range("K64").Value = 6314.25
$ws.Range("M64").Value = -6044.25
$ws.Range("H66").Value = 3725
$ws.Range("I66").Value = 2450
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 22050
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -18306
$ws.Range("N66").Value = -52488
$ws.Range("H67").Value = 3552.375
$ws.Range("I67").Value = 2104.75
$ws.Range("K67").Value = 6314.25
$ws.Range("M67").Value = -5378.25
$ws.Range("H82").Value = 5010000
$ws.Range("H85").Value = 5010000
$ws.Range("H122").Value = 738.5
$ws.Range("J122").Value = 1041
$ws.Range("L122").Value = 9369
$ws.Range("N122").Value = -14269

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13112.272
$ws.Range("I70").Value = 10970.786
$ws.Range("J70").Value = 16859.875
$ws.Range("K70").Value = 10970.786
$ws.Range("L70").Value = 16859.875
$ws.Range("M70").Value = -10700.786
$ws.Range("N70").Value = -17399.875
$ws.Range("H73").Value = 13112.272
$ws.Range("I73").Value = 10970.786
$ws.Range("J73").Value = 16859.875
$ws.Range("K73").Value = 10970.786
$ws.Range("L73").Value = 16859.875
$ws.Range("M73").Value = -10034.786
$ws.Range("N73").Value = -18731.875
$ws.Range("H107").Value = 693
$ws.Range("I107").Value = 693
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 693
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5356.607
$ws.Range("I46").Value = 6647.263
$ws.Range("J46").Value = 2631.889
$ws.Range("K46").Value = 6647.263
$ws.Range("L46").Value = 2631.889
$ws.Range("M46").Value = -6459.263
$ws.Range("N46").Value = -3007.889
$ws.Range("H132").Value = 10685.875
$ws.Range("I132").Value = 12848.333
$ws.Range("J132").Value = 4198.5
$ws.Range("K132").Value = 38544.999
$ws.Range("L132").Value = 12595.5
$ws.Range("M132").Value = -36014.999
$ws.Range("N132").Value = -17655.5
$ws.Range("H136").Value = 4642.6665
$ws.Range("I136").Value = 1714
$ws.Range("J136").Value = 10500
$ws.Range("K136").Value = 5142
$ws.Range("L136").Value = 31500
$ws.Range("M136").Value = -2592
$ws.Range("N136").Value = -36600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H46").Value = 74197.8
$ws.Range("J46").Value = 86497.25
$ws.Range("L46").Value = 86497.25
$ws.Range("N46").Value = -86959.25
$ws.Range("H62").Value = 8057.25
$ws.Range("I62").Value = 9729
$ws.Range("K62").Value = 9729
$ws.Range("M62").Value = -9105
$ws.Range("H65").Value = 8057.25
$ws.Range("I65").Value = 9729
$ws.Range("K65").Value = 48645
$ws.Range("M65").Value = -45525
$ws.Range("H107").Value = 1235
$ws.Range("J107").Value = 1778.75
$ws.Range("L107").Value = 5336.25
$ws.Range("N107").Value = -9176.25
$ws.Range("H122").Value = 2538.7827
$ws.Range("I122").Value = 1653.0588
$ws.Range("J122").Value = 5048.3335
$ws.Range("K122").Value = 4959.1764
$ws.Range("L122").Value = 15145.0005
$ws.Range("M122").Value = -2509.1764
$ws.Range("N122").Value = -20045.0005
$ws.Range("H126").Value = 4499.778
$ws.Range("I126").Value = 4499.778
$ws.Range("K126").Value = 13499.334
$ws.Range("M126").Value = -11029.334
$ws.Range("H134").Value = 74197.8
$ws.Range("J134").Value = 86497.25
$ws.Range("L134").Value = 259491.75
$ws.Range("N134").Value = -264561.75
$ws.Range("H136").Value = 6483.143
$ws.Range("I136").Value = 4426.5
$ws.Range("K136").Value = 13279.5
$ws.Range("M136").Value = -10729.5
